$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Insert a new paragraph right after the first (Heading1) paragraph that
#    contains a "Meta description" label (bold) followed by the description
#    text (regular).
# ---------------------------------------------------------------------------
$p1 = $d.Paragraphs(1)
$p1.Range.InsertParagraphAfter()
$metaPara = $d.Paragraphs(2)
$metaPara.Style = "Normal"

$metaStart = $metaPara.Range.Start
$metaEnd = $metaPara.Range.End
$metaFull = $d.Range($metaStart, $metaEnd - 1)
$metaFull.Text = "Meta description: Read our review of Don Juan's Peppers, a unique and engaging Mexican-themed slot game with special symbols and bonus features. Play for free!"

$boldLabel = "Meta description"
$boldRange = $d.Range($metaStart, $metaStart + $boldLabel.Length)
$boldRange.Font.Bold = 1

# ---------------------------------------------------------------------------
# 2) Remove the paragraph near the end of the document that duplicated the
#    title ("Play Don Juan's Peppers for Free - Slot Game Review") in bold.
# ---------------------------------------------------------------------------
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs($i)
    if ($para.Range.Text -eq "Play Don Juan's Peppers for Free - Slot Game Review`r") {
        $target = $para
    }
}
if ($target -ne $null) {
    $target.Range.Delete()
}

# ---------------------------------------------------------------------------
# 3) Replace the text of the final (italic) paragraph with the new image
#    generation prompt, keeping the italic formatting and straight quotes.
# ---------------------------------------------------------------------------
$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$lastStart = $lastPara.Range.Start
$lastEnd = $lastPara.Range.End
$lastText = $d.Range($lastStart, $lastEnd - 1)
$lastText.Text = "Please create a feature image fitting the game `"Don Juan's Peppers`". The image should be in a cartoon style and feature a happy Maya warrior with glasses."
